# Auto-generated Excel COM-interop script to apply the Zeromus_Profits.xlsx value updates.
# For each affected (sheet, row), set new numeric values in columns H:N and
# clear any cell that the target state no longer has (only column N in 3 specific rows).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2190.0576
$ws.Range("I17").Value = 350
$ws.Range("J17").Value = 2263.66
$ws.Range("K17").Value = 1050
$ws.Range("L17").Value = 6790.98
$ws.Range("M17").Value = -882
$ws.Range("N17").Value = -7126.98

$ws.Range("H132").Value = 3267.9832
$ws.Range("I132").Value = 2916.8958
$ws.Range("J132").Value = 4800
$ws.Range("K132").Value = 8750.687399999999
$ws.Range("L132").Value = 14400
$ws.Range("M132").Value = -6220.687399999999
$ws.Range("N132").Value = -19460

$ws.Range("H137").Value = 1199.194
$ws.Range("I137").Value = 799.8261
$ws.Range("J137").Value = 2074
$ws.Range("K137").Value = 2399.4783
$ws.Range("L137").Value = 6222
$ws.Range("M137").Value = 150.5217000000002
$ws.Range("N137").Value = -11322

$ws.Range("H138").Value = 1931.8704
$ws.Range("J138").Value = 3900.6667
$ws.Range("L138").Value = 11702.0001
$ws.Range("N138").Value = -21982.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 151
$ws.Range("I5").Value = 151
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 151
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -39
$ws.Range("N5").ClearContents()

$ws.Range("H45").Value = 2172.6
$ws.Range("I45").Value = 2532.2856
$ws.Range("K45").Value = 2532.2856
$ws.Range("M45").Value = -2155.2856

$ws.Range("H61").Value = 1211.4642
$ws.Range("I61").Value = 1124.2778
$ws.Range("J61").Value = 1368.4
$ws.Range("K61").Value = 1124.2778
$ws.Range("L61").Value = 1368.4
$ws.Range("M61").Value = -912.2778000000001
$ws.Range("N61").Value = -1792.4

$ws.Range("H74").Value = 7145610.5
$ws.Range("I74").Value = 8930880
$ws.Range("J74").Value = 4534.5713
$ws.Range("K74").Value = 8930880
$ws.Range("L74").Value = 4534.5713
$ws.Range("M74").Value = -8930006
$ws.Range("N74").Value = -6282.5713

$ws.Range("H77").Value = 7145610.5
$ws.Range("I77").Value = 8930880
$ws.Range("J77").Value = 4534.5713
$ws.Range("K77").Value = 44654400
$ws.Range("L77").Value = 22672.8565
$ws.Range("M77").Value = -44650032
$ws.Range("N77").Value = -31408.8565

$ws.Range("H110").Value = 6623.75
$ws.Range("I110").Value = 7737.278
$ws.Range("J110").Value = 3283.1667
$ws.Range("K110").Value = 7737.278
$ws.Range("L110").Value = 3283.1667
$ws.Range("M110").Value = -5692.278
$ws.Range("N110").Value = -7373.1667

$ws.Range("H123").Value = 52171.6
$ws.Range("J123").Value = 52171.6
$ws.Range("L123").Value = 52171.6
$ws.Range("N123").Value = -61971.6

$ws.Range("H136").Value = 1211.4642
$ws.Range("I136").Value = 1124.2778
$ws.Range("J136").Value = 1368.4
$ws.Range("K136").Value = 3372.8334
$ws.Range("L136").Value = 4105.200000000001
$ws.Range("M136").Value = -822.8334000000004
$ws.Range("N136").Value = -9205.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 151
$ws.Range("I4").Value = 151
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 151
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -36
$ws.Range("N4").ClearContents()

$ws.Range("H20").Value = 3741.8215
$ws.Range("I20").Value = 2341.3125
$ws.Range("J20").Value = 5609.1665
$ws.Range("K20").Value = 2341.3125
$ws.Range("L20").Value = 5609.1665
$ws.Range("M20").Value = -2094.3125
$ws.Range("N20").Value = -6103.1665

$ws.Range("H94").Value = 7070.6665
$ws.Range("I94").Value = 375.2
$ws.Range("J94").Value = 40548
$ws.Range("K94").Value = 375.2
$ws.Range("L94").Value = 40548
$ws.Range("M94").Value = 75.80000000000001
$ws.Range("N94").Value = -41450

$ws.Range("H134").Value = 2423.4546
$ws.Range("I134").Value = 2253.4666
$ws.Range("J134").Value = 2787.7144
$ws.Range("K134").Value = 6760.399800000001
$ws.Range("L134").Value = 8363.143199999999
$ws.Range("M134").Value = -4225.399800000001
$ws.Range("N134").Value = -13433.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3723596.5
$ws.Range("J31").Value = 2396.9412
$ws.Range("L31").Value = 2396.9412
$ws.Range("N31").Value = -2986.9412

$ws.Range("H34").Value = 3723596.5
$ws.Range("J34").Value = 2396.9412
$ws.Range("L34").Value = 2396.9412
$ws.Range("N34").Value = -2800.9412

$ws.Range("H58").Value = 1157.4286
$ws.Range("I58").Value = 728.6
$ws.Range("J58").Value = 1652.2307
$ws.Range("K58").Value = 728.6
$ws.Range("L58").Value = 1652.2307
$ws.Range("M58").Value = -525.6
$ws.Range("N58").Value = -2058.2307

$ws.Range("H132").Value = 1214.325
$ws.Range("I132").Value = 1035.9642
$ws.Range("J132").Value = 1630.5
$ws.Range("K132").Value = 3107.8926
$ws.Range("L132").Value = 4891.5
$ws.Range("M132").Value = -577.8925999999997
$ws.Range("N132").Value = -9951.5

$ws.Range("H134").Value = 1938.5111
$ws.Range("I134").Value = 2074.5264
$ws.Range("J134").Value = 1200.1428
$ws.Range("K134").Value = 6223.5792
$ws.Range("L134").Value = 3600.4284
$ws.Range("M134").Value = -3688.5792
$ws.Range("N134").Value = -8670.428400000001

$ws.Range("H136").Value = 1157.4286
$ws.Range("I136").Value = 728.6
$ws.Range("J136").Value = 1652.2307
$ws.Range("K136").Value = 2185.8
$ws.Range("L136").Value = 4956.6921
$ws.Range("M136").Value = 364.1999999999998
$ws.Range("N136").Value = -10056.6921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 375
$ws.Range("I107").Value = 368.25
$ws.Range("J107").Value = 378.85715
$ws.Range("K107").Value = 1104.75
$ws.Range("L107").Value = 1136.57145
$ws.Range("M107").Value = 815.25
$ws.Range("N107").Value = -4976.571449999999

$ws.Range("H113").Value = 1143.5883
$ws.Range("I113").Value = 522.125
$ws.Range("J113").Value = 1696
$ws.Range("K113").Value = 1566.375
$ws.Range("L113").Value = 5088
$ws.Range("M113").Value = 603.625
$ws.Range("N113").Value = -9428

$ws.Range("H131").Value = 1419328.2
$ws.Range("I131").Value = 4444735
$ws.Range("J131").Value = 1168.875
$ws.Range("K131").Value = 13334205
$ws.Range("L131").Value = 3506.625
$ws.Range("M131").Value = -13329165
$ws.Range("N131").Value = -13586.625

$ws.Range("H132").Value = 964.2857
$ws.Range("J132").Value = 964.2857
$ws.Range("L132").Value = 8678.5713
$ws.Range("N132").Value = -13738.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 53338616
$ws.Range("I70").Value = 100004640
$ws.Range("J70").Value = 6015.5713
$ws.Range("K70").Value = 100004640
$ws.Range("L70").Value = 6015.5713
$ws.Range("M70").Value = -100004370
$ws.Range("N70").Value = -6555.5713

$ws.Range("H73").Value = 53338616
$ws.Range("I73").Value = 100004640
$ws.Range("J73").Value = 6015.5713
$ws.Range("K73").Value = 100004640
$ws.Range("L73").Value = 6015.5713
$ws.Range("M73").Value = -100003704
$ws.Range("N73").Value = -7887.5713

$ws.Range("H102").Value = 1278.2069
$ws.Range("I102").Value = 1289.75
$ws.Range("J102").Value = 1222.8
$ws.Range("K102").Value = 1289.75
$ws.Range("L102").Value = 1222.8
$ws.Range("M102").Value = 332.25
$ws.Range("N102").Value = -4466.8

$ws.Range("H126").Value = 2516.1428
$ws.Range("I126").Value = 1551.9231
$ws.Range("K126").Value = 4655.7693
$ws.Range("M126").Value = -2185.7693

$ws.Range("H132").Value = 1784.4
$ws.Range("I132").Value = 1492.4
$ws.Range("J132").Value = 2368.4
$ws.Range("K132").Value = 4477.200000000001
$ws.Range("L132").Value = 7105.200000000001
$ws.Range("M132").Value = -1947.200000000001
$ws.Range("N132").Value = -12165.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4100
$ws.Range("I122").Value = 4225
$ws.Range("J122").Value = 3850
$ws.Range("K122").Value = 12675
$ws.Range("L122").Value = 11550
$ws.Range("M122").Value = -10225
$ws.Range("N122").Value = -16450

$ws.Range("H132").Value = 17865570
$ws.Range("I132").Value = 39077510
$ws.Range("J132").Value = 2881.7368
$ws.Range("K132").Value = 117232530
$ws.Range("L132").Value = 8645.2104
$ws.Range("M132").Value = -117230000
$ws.Range("N132").Value = -13705.2104

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 948.75
$ws.Range("I107").Value = 948.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2846.25
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -926.25
$ws.Range("N107").ClearContents()

$ws.Range("H123").Value = 48720.715
$ws.Range("J123").Value = 48720.715
$ws.Range("L123").Value = 48720.715
$ws.Range("N123").Value = -58520.715

$ws.Range("H132").Value = 4896.3076
$ws.Range("I132").Value = 2214.7
$ws.Range("J132").Value = 13835
$ws.Range("K132").Value = 6644.099999999999
$ws.Range("L132").Value = 41505
$ws.Range("M132").Value = -4114.099999999999
$ws.Range("N132").Value = -46565

$ws.Range("H136").Value = 933.46295
$ws.Range("I136").Value = 498.02127
$ws.Range("J136").Value = 3857.1428
$ws.Range("K136").Value = 1494.06381
$ws.Range("L136").Value = 11571.4284
$ws.Range("M136").Value = 1055.93619
$ws.Range("N136").Value = -16671.4284
